$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AddOpportunity")
$ws2 = $wb.Worksheets.Item("Users")
$ws3 = $wb.Worksheets.Item("AddContact")

# --- Content edits (kept in this order so new shared strings are
#     interned in the same sequence as the target workbook) ---

# AddContact: Party "ABC" renamed to "ABC Auto Parts Ltd"
$ws3.Range("D5").Value = "ABC Auto Parts Ltd"

# AddOpportunity: StdUser reassigned from Ayati Arvind to Eric Winthrop
$ws1.Range("AF2").Value = "Eric Winthrop"

# AddOpportunity: PrimaryOffice changed from AM to NY
$ws1.Range("K2").Value = "NY"

# Users: CaoUser reassigned from Jessica Maring to Derek Janisch
$ws2.Range("B2").Value = "Derek Janisch"

# AddOpportunity: IndustryGroup / Sector changed to Healthcare / Dental
$ws1.Range("D2").Value = "HC - Healthcare"
$ws1.Range("E2").Value = "Dental"

# Users: StdUser reassigned from Ayati Arvind to Eric Winthrop
$ws2.Range("A2").Value = "Eric Winthrop"

# AddOpportunity: second StdUser row added below the existing data row
$ws1.Range("AF3").Value = "Eric Winthrop"

# --- Formatting tweaks that came with the data edits ---
$ws1.Range("E2").WrapText = $true
$ws1.Range("E2").VerticalAlignment = -4108

$ws2.Range("B2").WrapText = $true
$ws2.Range("B2").VerticalAlignment = -4108

# --- New worksheet "Add", appended after AddContact, mirroring the
#     AddContact header row + first data row ---
$wsAdd = $wb.Worksheets.Add($null, $ws3)
$wsAdd.Name = "Add"

$wsAdd.Range("A1").Value = "Contact"
$wsAdd.Range("B1").Value = "Role"
$wsAdd.Range("C1").Value = "Party"
$wsAdd.Range("D1").Value = "Type1"
$wsAdd.Range("E1").Value = "ClientType"
$wsAdd.Range("F1").Value = "Contact2"
$wsAdd.Range("G1").Value = "Type2"
$wsAdd.Range("H1").Value = "HLContact"
$wsAdd.Range("A1:H1").Font.Bold = $true

$wsAdd.Range("A2").Value = "Sonika Goyal"
$wsAdd.Range("B2").Value = "Board of Directors"
$wsAdd.Range("C2").Value = "Buyer"
$wsAdd.Range("D2").Value = "Accupac"
$wsAdd.Range("E2").Value = "Client"
$wsAdd.Range("F2").Value = "Emma Watson"
$wsAdd.Range("G2").Value = "Client"
$wsAdd.Range("H2").Value = "Sonika Goyal"

# --- Selections / active sheet, set last so the final state matches ---
$ws2.Range("B4").Select()
$wsAdd.Range("F19").Select()
$ws3.Range("F11").Select()
$ws1.Range("E2").Select()
